# Add Group Name + Team Members name to Project Description Document
#
# The document currently contains a single paragraph with the title
# ("Betrayal at the House on the Hill by Nobody"). This adds a new
# paragraph right after it containing the team member names, inheriting
# the same paragraph/run formatting (Arial, size 11/22 half-points,
# line spacing 276, etc.) as the title paragraph.

$d = $word.ActiveDocument

# The last paragraph in the document (the title line).
$titlePara = $d.Paragraphs($d.Paragraphs.Count)

# Collapse a range to the very end of that paragraph (before its
# paragraph mark) and insert a new paragraph break there. The new
# paragraph inherits the formatting (pPr/rPr) of the paragraph it was
# split from, matching the Arial / 276 line-spacing formatting used
# throughout this document.
$endRange = $titlePara.Range
$endRange.SetRange($endRange.End, $endRange.End)
$endRange.InsertParagraphAfter()

# Fill in the text of the newly created (now last) paragraph.
$newPara = $d.Paragraphs($d.Paragraphs.Count)
$newRange = $newPara.Range
$newRange.SetRange($newRange.End, $newRange.End)
$newRange.Text = "(Caleb Corlett, Chris Anderson, Ethan Wyman, Patrick Storer, Ryan Nodarse)"
